$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 1.044366416687426
$ws.Cells.Item(2, 4).Value = 1.051754987288968
$ws.Cells.Item(2, 5).Value = 1.051835858366668
$ws.Cells.Item(2, 6).Value = 1.061815470365856
$ws.Cells.Item(2, 10).Value = 1.049431886253281
$ws.Cells.Item(2, 11).Value = 1.054505611127021
$ws.Cells.Item(2, 12).Value = 1.054586258161316
$ws.Cells.Item(2, 13).Value = 1.064538506498105
$ws.Cells.Item(2, 14).Value = 1.050922199903281

$ws.Cells.Item(3, 3).Value = 1.045742149378773
$ws.Cells.Item(3, 4).Value = 1.053009638720441
$ws.Cells.Item(3, 5).Value = 1.053073294039903
$ws.Cells.Item(3, 6).Value = 1.063166637302956
$ws.Cells.Item(3, 10).Value = 1.050452735914317
$ws.Cells.Item(3, 11).Value = 1.055571836503854
$ws.Cells.Item(3, 12).Value = 1.055635328350023
$ws.Cells.Item(3, 13).Value = 1.06570301774435
$ws.Cells.Item(3, 14).Value = 1.051944499287929

$ws.Cells.Item(4, 3).Value = 1.046632404330529
$ws.Cells.Item(4, 4).Value = 1.053821801990219
$ws.Cells.Item(4, 5).Value = 1.053874310918591
$ws.Cells.Item(4, 6).Value = 1.064041456525528
$ws.Cells.Item(4, 10).Value = 1.051112912620487
$ws.Cells.Item(4, 11).Value = 1.05626151719647
$ws.Cells.Item(4, 12).Value = 1.056313898129421
$ws.Cells.Item(4, 13).Value = 1.066456516076862
$ws.Cells.Item(4, 14).Value = 1.052605613520744

$ws.Cells.Item(5, 3).Value = 1.047006686586925
$ws.Cells.Item(5, 4).Value = 1.054163314693475
$ws.Cells.Item(5, 5).Value = 1.054211135995033
$ws.Cells.Item(5, 6).Value = 1.064409358877251
$ws.Cells.Item(5, 10).Value = 1.051390362047624
$ws.Cells.Item(5, 11).Value = 1.056551404318079
$ws.Cells.Item(5, 12).Value = 1.056599111670254
$ws.Cells.Item(5, 13).Value = 1.066773285011319
$ws.Cells.Item(5, 14).Value = 1.052883456957891

$ws.Cells.Item(6, 3).Value = 1.047069531452981
$ws.Cells.Item(6, 4).Value = 1.054220660905103
$ws.Cells.Item(6, 5).Value = 1.054267695032027
$ws.Cells.Item(6, 6).Value = 1.064471138912841
$ws.Cells.Item(6, 10).Value = 1.051436941886827
$ws.Cells.Item(6, 11).Value = 1.05660007450799
$ws.Cells.Item(6, 12).Value = 1.056646996993569
$ws.Cells.Item(6, 13).Value = 1.066826471881926
$ws.Cells.Item(6, 14).Value = 1.05293010294581

$ws.Cells.Item(7, 3).Value = 1.046637405428698
$ws.Cells.Item(7, 4).Value = 1.053826364984879
$ws.Cells.Item(7, 5).Value = 1.053878811283686
$ws.Cells.Item(7, 6).Value = 1.064046371946768
$ws.Cells.Item(7, 10).Value = 1.051116620259461
$ws.Cells.Item(7, 11).Value = 1.056265390894049
$ws.Cells.Item(7, 12).Value = 1.056317709388507
$ws.Cells.Item(7, 13).Value = 1.06646074876254
$ws.Cells.Item(7, 14).Value = 1.052609326424991

$ws.Cells.Item(8, 3).Value = 1.044831341031003
$ws.Cells.Item(8, 4).Value = 1.052178937991249
$ws.Cells.Item(8, 5).Value = 1.052253992242456
$ws.Cells.Item(8, 6).Value = 1.062271996856341
$ws.Cells.Item(8, 10).Value = 1.049776967580925
$ws.Cells.Item(8, 11).Value = 1.054865997819324
$ws.Cells.Item(8, 12).Value = 1.054940849278645
$ws.Cells.Item(8, 13).Value = 1.064932063894522
$ws.Cells.Item(8, 14).Value = 1.051267771285991

$ws.Cells.Item(9, 3).Value = 1.041649145447722
$ws.Cells.Item(9, 4).Value = 1.049278291101655
$ws.Cells.Item(9, 5).Value = 1.049393136725764
$ws.Cells.Item(9, 6).Value = 1.05914920193163
$ws.Cells.Item(9, 10).Value = 1.047413299785727
$ws.Cells.Item(9, 11).Value = 1.052398153224373
$ws.Cells.Item(9, 12).Value = 1.052512634553902
$ws.Cells.Item(9, 13).Value = 1.062238068275244
$ws.Cells.Item(9, 14).Value = 1.048900746811408

$ws.Cells.Item(10, 3).Value = 1.039527660251702
$ws.Cells.Item(10, 4).Value = 1.04734592356721
$ws.Cells.Item(10, 5).Value = 1.047487271135892
$ws.Cells.Item(10, 6).Value = 1.057069762400185
$ws.Cells.Item(10, 10).Value = 1.045835332076312
$ws.Cells.Item(10, 11).Value = 1.050751468326262
$ws.Cells.Item(10, 12).Value = 1.050892322931896
$ws.Cells.Item(10, 13).Value = 1.060441732524757
$ws.Cells.Item(10, 14).Value = 1.047320538206849

$ws.Cells.Item(11, 3).Value = 1.038608975565838
$ws.Cells.Item(11, 4).Value = 1.046509480238738
$ws.Cells.Item(11, 5).Value = 1.046662298636129
$ws.Cells.Item(11, 6).Value = 1.056169875254395
$ws.Cells.Item(11, 10).Value = 1.04515150093701
$ws.Cells.Item(11, 11).Value = 1.050038055707306
$ws.Cells.Item(11, 12).Value = 1.050190320172058
$ws.Cells.Item(11, 13).Value = 1.0596637816879
$ws.Cells.Item(11, 14).Value = 1.046635735948893

$ws.Cells.Item(12, 3).Value = 1.038267720380208
$ws.Cells.Item(12, 4).Value = 1.046198827334766
$ws.Cells.Item(12, 5).Value = 1.046355905913463
$ws.Cells.Item(12, 6).Value = 1.055835692144694
$ws.Cells.Item(12, 10).Value = 1.04489740832854
$ws.Cells.Item(12, 11).Value = 1.049773001264149
$ws.Cells.Item(12, 12).Value = 1.04992950245072
$ws.Cells.Item(12, 13).Value = 1.059374793956056
$ws.Cells.Item(12, 14).Value = 1.046381282499774

$ws.Cells.Item(13, 3).Value = 1.038340921519681
$ws.Cells.Item(13, 4).Value = 1.046265461708766
$ws.Cells.Item(13, 5).Value = 1.046421626489293
$ws.Cells.Item(13, 6).Value = 1.05590737221921
$ws.Cells.Item(13, 10).Value = 1.044951916053832
$ws.Cells.Item(13, 11).Value = 1.04982985916205
$ws.Cells.Item(13, 12).Value = 1.049985451619737
$ws.Cells.Item(13, 13).Value = 1.059436783860785
$ws.Cells.Item(13, 14).Value = 1.046435867632289

$ws.Cells.Item(14, 3).Value = 1.038580767634264
$ws.Cells.Item(14, 4).Value = 1.046483800788703
$ws.Cells.Item(14, 5).Value = 1.046636971345521
$ws.Cells.Item(14, 6).Value = 1.056142250064665
$ws.Cells.Item(14, 10).Value = 1.045130499360715
$ws.Cells.Item(14, 11).Value = 1.050016147491192
$ws.Cells.Item(14, 12).Value = 1.050168762192914
$ws.Cells.Item(14, 13).Value = 1.059639894342761
$ws.Cells.Item(14, 14).Value = 1.046614704547951

$ws.Cells.Item(15, 3).Value = 1.038728542642932
$ws.Cells.Item(15, 4).Value = 1.046618331820897
$ws.Cells.Item(15, 5).Value = 1.046769657462061
$ws.Cells.Item(15, 6).Value = 1.056286975891406
$ws.Cells.Item(15, 10).Value = 1.045240518774252
$ws.Cells.Item(15, 11).Value = 1.05013091767357
$ws.Cells.Item(15, 12).Value = 1.050281697497581
$ws.Cells.Item(15, 13).Value = 1.059765034380407
$ws.Cells.Item(15, 14).Value = 1.046724880201675

$ws.Cells.Item(16, 3).Value = 1.039588628509275
$ws.Cells.Item(16, 4).Value = 1.047401441264886
$ws.Cells.Item(16, 5).Value = 1.047542027462982
$ws.Cells.Item(16, 6).Value = 1.057129495671544
$ws.Cells.Item(16, 10).Value = 1.04588070362447
$ws.Cells.Item(16, 11).Value = 1.050798806753473
$ws.Cells.Item(16, 12).Value = 1.050938903922175
$ws.Cells.Item(16, 13).Value = 1.060493359633398
$ws.Cells.Item(16, 14).Value = 1.04736597418781

$ws.Cells.Item(17, 3).Value = 1.040128116284321
$ws.Cells.Item(17, 4).Value = 1.047892738944828
$ws.Cells.Item(17, 5).Value = 1.04802658739586
$ws.Cells.Item(17, 6).Value = 1.057658123746517
$ws.Cells.Item(17, 10).Value = 1.046282122408862
$ws.Cells.Item(17, 11).Value = 1.051217650106969
$ws.Cells.Item(17, 12).Value = 1.051351043726761
$ws.Cells.Item(17, 13).Value = 1.0609501828799
$ws.Cells.Item(17, 14).Value = 1.047767963032919

$ws.Cells.Item(18, 3).Value = 1.040442784168372
$ws.Cells.Item(18, 4).Value = 1.048179332254409
$ws.Cells.Item(18, 5).Value = 1.048309250231855
$ws.Cells.Item(18, 6).Value = 1.057966514356165
$ws.Cells.Item(18, 10).Value = 1.046516209300371
$ws.Cells.Item(18, 11).Value = 1.051461917198398
$ws.Cells.Item(18, 12).Value = 1.051591399786307
$ws.Cells.Item(18, 13).Value = 1.06121662825493
$ws.Cells.Item(18, 14).Value = 1.048002382354664

$ws.Cells.Item(19, 3).Value = 1.040550076871773
$ws.Cells.Item(19, 4).Value = 1.048277057963749
$ws.Cells.Item(19, 5).Value = 1.04840563567
$ws.Cells.Item(19, 6).Value = 1.058071676310927
$ws.Cells.Item(19, 10).Value = 1.04659601783105
$ws.Cells.Item(19, 11).Value = 1.051545199754903
$ws.Cells.Item(19, 12).Value = 1.051673348615695
$ws.Cells.Item(19, 13).Value = 1.061307477323297
$ws.Cells.Item(19, 14).Value = 1.048082304222612

$ws.Cells.Item(20, 3).Value = 1.040070235028806
$ws.Cells.Item(20, 4).Value = 1.047840024515357
$ws.Cells.Item(20, 5).Value = 1.047974595909623
$ws.Cells.Item(20, 6).Value = 1.057601401730623
$ws.Cells.Item(20, 10).Value = 1.046239059569863
$ws.Cells.Item(20, 11).Value = 1.05117271604495
$ws.Cells.Item(20, 12).Value = 1.051306828991294
$ws.Cells.Item(20, 13).Value = 1.06090117132428
$ws.Cells.Item(20, 14).Value = 1.047724839039749

$ws.Cells.Item(21, 3).Value = 1.038510139377034
$ws.Cells.Item(21, 4).Value = 1.046419504311286
$ws.Cells.Item(21, 5).Value = 1.046573556607911
$ws.Cells.Item(21, 6).Value = 1.056073082351266
$ws.Cells.Item(21, 10).Value = 1.045077913461565
$ws.Cells.Item(21, 11).Value = 1.04996129193629
$ws.Cells.Item(21, 12).Value = 1.050114783551098
$ws.Cells.Item(21, 13).Value = 1.059580084006143
$ws.Cells.Item(21, 14).Value = 1.046562043970793

$ws.Cells.Item(22, 3).Value = 1.037529157455006
$ws.Cells.Item(22, 4).Value = 1.045526593110581
$ws.Cells.Item(22, 5).Value = 1.045692890434128
$ws.Cells.Item(22, 6).Value = 1.055112599656356
$ws.Cells.Item(22, 10).Value = 1.044347348358119
$ws.Cells.Item(22, 11).Value = 1.049199265728368
$ws.Cells.Item(22, 12).Value = 1.0493649333783
$ws.Cells.Item(22, 13).Value = 1.058749334174753
$ws.Cells.Item(22, 14).Value = 1.045830441381103

$ws.Cells.Item(23, 3).Value = 1.038049204410176
$ws.Cells.Item(23, 4).Value = 1.045999921912478
$ws.Cells.Item(23, 5).Value = 1.046159728200397
$ws.Cells.Item(23, 6).Value = 1.055621729948876
$ws.Cells.Item(23, 10).Value = 1.044734683763612
$ws.Cells.Item(23, 11).Value = 1.049603264987206
$ws.Cells.Item(23, 12).Value = 1.049762478629855
$ws.Cells.Item(23, 13).Value = 1.059189743779694
$ws.Cells.Item(23, 14).Value = 1.046218326847299

$ws.Cells.Item(24, 3).Value = 1.040096389095662
$ws.Cells.Item(24, 4).Value = 1.047863843811324
$ws.Cells.Item(24, 5).Value = 1.047998088540091
$ws.Cells.Item(24, 6).Value = 1.057627031810961
$ws.Cells.Item(24, 10).Value = 1.046258517979599
$ws.Cells.Item(24, 11).Value = 1.051193019928095
$ws.Cells.Item(24, 12).Value = 1.051326807844723
$ws.Cells.Item(24, 13).Value = 1.060923317573232
$ws.Cells.Item(24, 14).Value = 1.04774432508266

$ws.Cells.Item(25, 3).Value = 1.042471805672415
$ws.Cells.Item(25, 4).Value = 1.050027919167129
$ws.Cells.Item(25, 5).Value = 1.050132482710053
$ws.Cells.Item(25, 6).Value = 1.059956078459365
$ws.Cells.Item(25, 10).Value = 1.048024739200765
$ws.Cells.Item(25, 11).Value = 1.053036395380807
$ws.Cells.Item(25, 12).Value = 1.053140639990817
$ws.Cells.Item(25, 13).Value = 1.062934578261756
$ws.Cells.Item(25, 14).Value = 1.049513054540548
